$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 585.2857
$ws.Range("I5").Value = 381.72726
$ws.Range("J5").Value = 1331.6666
$ws.Range("K5").Value = 381.72726
$ws.Range("L5").Value = 1331.6666
$ws.Range("M5").Value = -266.72726
$ws.Range("N5").Value = -1561.6666

$ws.Range("H15").Value = 3330.2778
$ws.Range("I15").Value = 3330.2778
$ws.Range("K15").Value = 9990.8334
$ws.Range("M15").Value = -9821.8334

$ws.Range("H33").Value = 1405
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").Value = $null

$ws.Range("H98").Value = 1046.8
$ws.Range("I98").Value = 1058.5
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 1058.5
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = 439.5
$ws.Range("N98").Value = -3996

$ws.Range("H112").Value = 2721
$ws.Range("J112").Value = 2853.3845
$ws.Range("L112").Value = 8560.1535
$ws.Range("N112").Value = -10776.1535

$ws.Range("H115").Value = 734.9091
$ws.Range("I115").Value = 734.9091
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 2204.7273
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -637.7273
$ws.Range("N115").Value = $null

$ws.Range("H122").Value = 1046.8
$ws.Range("I122").Value = 1058.5
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 3175.5
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -725.5
$ws.Range("N122").Value = -7900

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2316.875
$ws.Range("I2").Value = 2719.3333
$ws.Range("K2").Value = 2719.3333
$ws.Range("M2").Value = -2606.3333

$ws.Range("H4").Value = 34.8
$ws.Range("I4").Value = 41.5
$ws.Range("J4").Value = 8
$ws.Range("K4").Value = 41.5
$ws.Range("L4").Value = 8
$ws.Range("M4").Value = 74.5
$ws.Range("N4").Value = -240

$ws.Range("H45").Value = 2301.7693
$ws.Range("I45").Value = 2406.1667
$ws.Range("K45").Value = 2406.1667
$ws.Range("M45").Value = -2029.1667

$ws.Range("H63").Value = 8601
$ws.Range("I63").Value = 1000
$ws.Range("K63").Value = 1000
$ws.Range("M63").Value = -314

$ws.Range("H66").Value = 8601
$ws.Range("I66").Value = 1000
$ws.Range("K66").Value = 5000
$ws.Range("M66").Value = -1568

$ws.Range("H74").Value = 7356
$ws.Range("I74").Value = 7091.7
$ws.Range("K74").Value = 7091.7
$ws.Range("M74").Value = -6217.7

$ws.Range("H77").Value = 7356
$ws.Range("I77").Value = 7091.7
$ws.Range("K77").Value = 35458.5
$ws.Range("M77").Value = -31090.5

$ws.Range("H116").Value = 2316.875
$ws.Range("I116").Value = 2719.3333
$ws.Range("K116").Value = 2719.3333
$ws.Range("M116").Value = -425.3332999999998

$ws.Range("H122").Value = 3869.5454
$ws.Range("I122").Value = 3507.111
$ws.Range("K122").Value = 10521.333
$ws.Range("M122").Value = -8071.332999999999

$ws.Range("H132").Value = 2032.5
$ws.Range("I132").Value = 1276.8667
$ws.Range("J132").Value = 4299.4
$ws.Range("K132").Value = 3830.6001
$ws.Range("L132").Value = 12898.2
$ws.Range("M132").Value = -1300.6001
$ws.Range("N132").Value = -17958.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2316.875
$ws.Range("I3").Value = 2719.3333
$ws.Range("K3").Value = 2719.3333
$ws.Range("M3").Value = -2605.3333

$ws.Range("H105").Value = 3664
$ws.Range("I105").Value = 3596.8
$ws.Range("K105").Value = 3596.8
$ws.Range("M105").Value = -1849.8

$ws.Range("H107").Value = 1372.125
$ws.Range("I107").Value = 1372.125
$ws.Range("K107").Value = 1372.125
$ws.Range("M107").Value = 547.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 6667083.5
$ws.Range("I6").Value = 20000000
$ws.Range("K6").Value = 20000000
$ws.Range("M6").Value = -19999887

$ws.Range("H16").Value = 856.5
$ws.Range("I16").Value = 959.2
$ws.Range("J16").Value = 343
$ws.Range("K16").Value = 959.2
$ws.Range("L16").Value = 343
$ws.Range("M16").Value = -672.2
$ws.Range("N16").Value = -917

$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = $null
$ws.Range("N94").Value = $null

$ws.Range("H113").Value = 856.5
$ws.Range("I113").Value = 959.2
$ws.Range("J113").Value = 343
$ws.Range("K113").Value = 959.2
$ws.Range("L113").Value = 343
$ws.Range("M113").Value = 1210.8
$ws.Range("N113").Value = -4683

$ws.Range("H134").Value = 5193.1875
$ws.Range("I134").Value = 4651.364
$ws.Range("J134").Value = 6385.2
$ws.Range("K134").Value = 13954.092
$ws.Range("L134").Value = 19155.6
$ws.Range("M134").Value = -11419.092
$ws.Range("N134").Value = -24225.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3851.6667
$ws.Range("I132").Value = 3615.25
$ws.Range("K132").Value = 10845.75
$ws.Range("M132").Value = -8315.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2338.0833
$ws.Range("I55").Value = 6761.3335
$ws.Range("J55").Value = 863.6667
$ws.Range("K55").Value = 6761.3335
$ws.Range("L55").Value = 863.6667
$ws.Range("M55").Value = -6588.3335
$ws.Range("N55").Value = -1209.6667

$ws.Range("H122").Value = 4996.5
$ws.Range("I122").Value = 4996.5
$ws.Range("K122").Value = 14989.5
$ws.Range("M122").Value = -12539.5

$ws.Range("H136").Value = 3055.4285
$ws.Range("I136").Value = 3231.3333
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 9693.999899999999
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -7143.999899999999
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1596.75
$ws.Range("I126").Value = 1559.5555
$ws.Range("J126").Value = 1627.1818
$ws.Range("K126").Value = 4678.666499999999
$ws.Range("L126").Value = 4881.5454
$ws.Range("M126").Value = -2208.666499999999
$ws.Range("N126").Value = -9821.545399999999

$ws.Range("H132").Value = 2736.8
$ws.Range("I132").Value = 2413
$ws.Range("J132").Value = 3492.3333
$ws.Range("K132").Value = 7239
$ws.Range("L132").Value = 10476.9999
$ws.Range("M132").Value = -4709
$ws.Range("N132").Value = -15536.9999

$ws.Range("H136").Value = 5543.143
$ws.Range("I136").Value = 5543.143
$ws.Range("K136").Value = 16629.429
$ws.Range("M136").Value = -14079.429
